# Media.xlsx — add new event rows, fix the "10AM" -> "10am" typo, and
# hyperlink the Website column for the rows that have a clickable URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix the existing row 2 start time ("10AM" -> "10am") ---------------
$ws.Range("B2").Value = "10am"

# --- new data rows (3-17) ------------------------------------------------
# Columns: A=Event Date  B=Start Time  C=Venue  D=Post code
#          E=Cover  F=Media Type  G=Website
$rows = @(
    @{ Row=3;  Date=[DateTime]"2025-08-03"; Time="10am"; Venue="Maverick Gym";         PostCode=$null;     Cover="Alex Hulme";     Media="Photography"; Web="https://jagershoots.com" },
    @{ Row=4;  Date=[DateTime]"2025-08-10"; Time="10am"; Venue="349 Barbell";          PostCode="SP4 6AT"; Cover="Labibur Rahman"; Media="Photography"; Web="https://thephotolabx.myshopify.com/" },
    @{ Row=5;  Date=[DateTime]"2025-08-17"; Time="9am";  Venue="Raw Strength Gym";     PostCode=$null;     Cover="Alex Hulme";     Media="Photography"; Web="https://jagershoots.com" },
    @{ Row=6;  Date=[DateTime]"2025-08-24"; Time="10am"; Venue="Stag Fitness Centre";  PostCode=$null;     Cover="Mike Melladay";  Media="Photography"; Web="https://melladaymedia.co.uk" },
    @{ Row=7;  Date=[DateTime]"2025-09-07"; Time="10am"; Venue="Lincoln Lifting";      PostCode=$null;     Cover="Alex Hulme";     Media="Photography"; Web="https://jagershoots.com" },
    @{ Row=8;  Date=[DateTime]"2025-09-28"; Time="10am"; Venue="Spartan Fitness Gym";  PostCode=$null;     Cover="Alex Hulme";     Media="Photography"; Web="https://jagershoots.com" },
    @{ Row=9;  Date=[DateTime]"2025-10-05"; Time="9am";  Venue="349 Barbell";          PostCode="SP4 6AT"; Cover="Alex Hulme";     Media="Photography"; Web="https://jagershoots.com" },
    @{ Row=10; Date=[DateTime]"2025-10-05"; Time="9am";  Venue="349 Barbell";          PostCode="SP4 6AT"; Cover="Sam Taylor";     Media="Videography"; Web=$null },
    @{ Row=11; Date=[DateTime]"2025-11-01"; Time="9am";  Venue="Nottingham Strong";    PostCode="NG7 2FH"; Cover=$null;            Media=$null;         Web=$null },
    @{ Row=12; Date=[DateTime]"2025-11-02"; Time="9am";  Venue="Nottingham Strong";    PostCode="NG7 2FH"; Cover=$null;            Media=$null;         Web=$null },
    @{ Row=13; Date=[DateTime]"2025-11-16"; Time="10am"; Venue="Iron Warehouse Gym";   PostCode=$null;     Cover=$null;            Media=$null;         Web=$null },
    @{ Row=14; Date=[DateTime]"2025-11-23"; Time="10am"; Venue="Maverick Gym";         PostCode=$null;     Cover="Mike Melladay";  Media="Photography"; Web="https://melladaymedia.co.uk" },
    @{ Row=15; Date=[DateTime]"2025-11-30"; Time="10am"; Venue="349 Barbell";          PostCode="SP4 6AT"; Cover=$null;            Media=$null;         Web=$null },
    @{ Row=16; Date=[DateTime]"2025-12-06"; Time="9am";  Venue="Nottingham Strong";    PostCode="NG7 2FH"; Cover=$null;            Media=$null;         Web=$null },
    @{ Row=17; Date=[DateTime]"2025-12-07"; Time="9am";  Venue="Nottingham Strong";    PostCode="NG7 2FH"; Cover=$null;            Media=$null;         Web=$null }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Date
    $ws.Range("B$n").Value = $r.Time
    $ws.Range("C$n").Value = $r.Venue
    if ($r.PostCode) { $ws.Range("D$n").Value = $r.PostCode }
    if ($r.Cover)    { $ws.Range("E$n").Value = $r.Cover }
    if ($r.Media)    { $ws.Range("F$n").Value = $r.Media }
    if ($r.Web)      { $ws.Range("G$n").Value = $r.Web }
}

# --- turn the Website cells that should be clickable into hyperlinks ----
$hyperlinkRows = @(3, 5, 6, 7, 8, 9, 14)
foreach ($n in $hyperlinkRows) {
    $target = ($rows | Where-Object { $_.Row -eq $n }).Web
    $ws.Hyperlinks.Add($ws.Range("G$n"), $target)
}

# --- match the final selection left by the author -----------------------
[void]$ws.Range("C17:D17").Select()

Write-Host "done"
